$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set/update cell values
$ws.Range("AB2").Value = 0.02025057688989229
$ws.Range("AC2").Value = -0.0201285154436341
$ws.Range("AD2").Value = 2921.6
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2921.6
$ws.Range("AG2").Value = 2817.847
$ws.Range("AH2").Value = 0.7641965943867542
$ws.Range("AI2").Value = 0.712029635406512
$ws.Range("AJ2").Value = 0.7576187432901528
$ws.Range("AK2").Value = 0.7045591553032207
$ws.Range("AL2").Value = 72.59
$ws.Range("AM2").Value = 67.521
$ws.Range("AN2").Value = 51.07692307692307
$ws.Range("AO2").Value = 0.6116544978647196
$ws.Range("AP2").Value = 49.26305944055943
$ws.Range("AQ2").Value = 0.6575731994490602
$ws.Range("D2").Value = 0.0626
$ws.Range("E2").Value = 0.0693
$ws.Range("G2").Value = 0.1605174750987211
$ws.Range("H2").Value = 0.1579536747804562
$ws.Range("I2").Value = 0.0654211115695173
$ws.Range("J2").Value = 0.05551870390092661
$ws.Range("K2").Value = 21.18
$ws.Range("L2").Value = 0.03120763835681028
$ws.Range("M2").Value = 5.01086
$ws.Range("N2").Value = 0.005558358291735996
$ws.Range("O2").Value = 0.2365845136921624
$ws.Range("P2").Value = 4.95986
$ws.Range("Q2").Value = 0.005501785912368275
$ws.Range("R2").Value = 0.234176581680831
$ws.Range("S2").Value = 0.051
$ws.Range("T2").Value = 0.01017789361506807
$ws.Range("U2").Value = 103.753
$ws.Range("V2").Value = 0.1150892956184138
$ws.Range("W2").Value = 0.06986505322434941
$ws.Range("X2").Value = 0.03911320748566795
$ws.Range("Y2").Value = 0.03075184573868146
$ws.Range("Z2").Value = 0.1895512750331242
$ws.Range("AA3").Value = 0.07642903018625562
$ws.Range("AB3").Value = 0.01884388895845857
$ws.Range("AC3").Value = 0.05758514122779705
$ws.Range("AD3").Value = 146
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 146
$ws.Range("AG3").Value = 126.5
$ws.Range("AH3").Value = 0.4515929477265697
$ws.Range("AI3").Value = 0.8138238573021181
$ws.Range("AJ3").Value = 0.4163923633969717
$ws.Range("AK3").Value = 0.791119449656035
$ws.Range("AL3").Value = 6.59
$ws.Range("AM3").Value = 6.547
$ws.Range("AN3").Value = 6.854460093896713
$ws.Range("AO3").Value = 1.805766312594841
$ws.Range("AP3").Value = 5.938967136150235
$ws.Range("AQ3").Value = 1.817626393768138
$ws.Range("B3").Value = 'Infront ASA (OB:INFRO)'
$ws.Range("G3").Value = 0.08827067669172932
$ws.Range("H3").Value = 0.07518796992481203
$ws.Range("I3").Value = 0.08947368421052632
$ws.Range("J3").Value = 0.08947368421052632
$ws.Range("K3").Value = -9.24
$ws.Range("L3").Value = -0.06947368421052631
$ws.Range("M3").Value = 0.102
$ws.Range("N3").Value = 0.0005752961082910321
$ws.Range("O3").Value = -0.01103896103896104
$ws.Range("P3").Value = 0.102
$ws.Range("Q3").Value = 0.0005752961082910321
$ws.Range("R3").Value = -0.01103896103896104
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 19.5
$ws.Range("V3").Value = 0.1099830795262267
$ws.Range("W3").Value = -0.2610169491525424
$ws.Range("X3").Value = 0.02202895262419433
$ws.Range("Y3").Value = -0.2830459017767367
$ws.Range("Z3").Value = 0.8542068079640333
$ws.Range("AA4").Value = 0.02877379371403276
$ws.Range("AB4").Value = 0.02430223063087752
$ws.Range("AC4").Value = 0.004471563083155241
$ws.Range("AD4").Value = 1092.8
$ws.Range("AF4").Value = 1092.8
$ws.Range("AG4").Value = 1054
$ws.Range("AH4").Value = 0.795226313491486
$ws.Range("AI4").Value = 0.7183802261372602
$ws.Range("AJ4").Value = 0.7892766212370824
$ws.Range("AK4").Value = 0.7110091743119266
$ws.Range("AL4").Value = 66
$ws.Range("AM4").Value = 65.994
$ws.Range("AN4").Value = 30.44011142061281
$ws.Range("AO4").Value = 0.4924242424242424
$ws.Range("AP4").Value = 29.35933147632312
$ws.Range("AQ4").Value = 0.4924690123344547
$ws.Range("B4").Value = 'Axactor SE (OB:AXA)'
$ws.Range("G4").Value = 0.3739899961523663
$ws.Range("H4").Value = 0.3739899961523663
$ws.Range("I4").Value = 0.1250480954213159
$ws.Range("J4").Value = 0.1250480954213159
$ws.Range("K4").Value = -14.9
$ws.Range("L4").Value = -0.05732974220854176
$ws.Range("U4").Value = 38.8
$ws.Range("V4").Value = 0.1378820184790334
$ws.Range("W4").Value = -0.050236008091706
$ws.Range("X4").Value = 0.04052813835448431
$ws.Range("Y4").Value = -0.0907641464461903
$ws.Range("Z4").Value = 0.2301018149623727
$ws.Range("AB5").Value = 0.02009872798824948
$ws.Range("AC5").Value = -0.02009872798824948
$ws.Range("AD5").Value = 65.40000000000001
$ws.Range("AF5").Value = 65.40000000000001
$ws.Range("AG5").Value = 65.05800000000001
$ws.Range("AH5").Value = 0.767605633802817
$ws.Range("AI5").Value = 0.5528317836010144
$ws.Range("AJ5").Value = 0.7666690235452167
$ws.Range("AK5").Value = 0.5515352922226556
$ws.Range("B5").Value = 'Sunndal Sparebank (OB:SUNSB)'
$ws.Range("K5").Value = 3.99
$ws.Range("L5").Value = 0.387378640776699
$ws.Range("M5").Value = 1.12608
$ws.Range("N5").Value = 0.05687272727272727
$ws.Range("O5").Value = 0.2822255639097744
$ws.Range("P5").Value = 1.12608
$ws.Range("Q5").Value = 0.05687272727272727
$ws.Range("R5").Value = 0.2822255639097744
$ws.Range("U5").Value = 0.342
$ws.Range("V5").Value = 0.01727272727272727
$ws.Range("W5").Value = 0.07702702702702703
$ws.Range("X5").Value = 0.03701925376761897
$ws.Range("Y5").Value = 0.04000777325940806
$ws.Range("Z5").Value = 0.07469722242367104
$ws.Range("AB6").Value = 0.02015830289901872
$ws.Range("AC6").Value = -0.02015830289901872
$ws.Range("AD6").Value = 73.8
$ws.Range("AF6").Value = 73.8
$ws.Range("AG6").Value = 73.32899999999999
$ws.Range("AH6").Value = 0.7826086956521739
$ws.Range("AI6").Value = 0.6589285714285714
$ws.Range("AJ6").Value = 0.7815174413027955
$ws.Range("AK6").Value = 0.6574881869289602
$ws.Range("B6").Value = 'Tysnes Sparebank (OB:TYSB)'
$ws.Range("K6").Value = 2.66
$ws.Range("L6").Value = 0.3604336043360434
$ws.Range("M6").Value = 0.72178
$ws.Range("N6").Value = 0.03520878048780487
$ws.Range("O6").Value = 0.2713458646616541
$ws.Range("P6").Value = 0.72178
$ws.Range("Q6").Value = 0.03520878048780487
$ws.Range("R6").Value = 0.2713458646616541
$ws.Range("U6").Value = 0.471
$ws.Range("V6").Value = 0.02297560975609756
$ws.Range("W6").Value = 0.07492957746478873
$ws.Range("X6").Value = 0.0388145933354861
$ws.Range("Y6").Value = 0.03611498412930263
$ws.Range("Z6").Value = 0.07942315970727508
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0.02034285088076586
$ws.Range("AC7").Value = -0.02034285088076586
$ws.Range("AD7").Value = 230.9
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 230.9
$ws.Range("AG7").Value = 229.76
$ws.Range("AH7").Value = 0.8290843806104129
$ws.Range("AI7").Value = 0.6690814256737178
$ws.Range("AJ7").Value = 0.8283818863570811
$ws.Range("AK7").Value = 0.6679846493778346
$ws.Range("AM7").Value = 0
$ws.Range("B7").Value = 'Melhus Sparebank (OB:MELG)'
$ws.Range("D7").Value = 0.0626
$ws.Range("E7").Value = 0.0693
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 9.27
$ws.Range("L7").Value = 0.4138392857142857
$ws.Range("M7").Value = 2.96
$ws.Range("N7").Value = 0.06218487394957983
$ws.Range("O7").Value = 0.3193096008629989
$ws.Range("P7").Value = 2.96
$ws.Range("Q7").Value = 0.06218487394957983
$ws.Range("R7").Value = 0.3193096008629989
$ws.Range("U7").Value = 1.14
$ws.Range("V7").Value = 0.02394957983193277
$ws.Range("W7").Value = 0.07909556313993174
$ws.Range("X7").Value = 0.04637658761120361
$ws.Range("Y7").Value = 0.03271897552872813
$ws.Range("Z7").Value = 0.05468135902706239
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0.0389648442663487
$ws.Range("AC8").Value = -0.0389648442663487
$ws.Range("AD8").Value = 1312.7
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 1312.7
$ws.Range("AG8").Value = 1269.2
$ws.Range("AH8").Value = 0.7871791796593908
$ws.Range("AI8").Value = 0.7184216287215411
$ws.Range("AJ8").Value = 0.7814789729696447
$ws.Range("AK8").Value = 0.711554633626731
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = -5.02
$ws.Range("AQ8").Value = -0
$ws.Range("B8").Value = 'B2Holding ASA (OB:B2H)'
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 29.4
$ws.Range("L8").Value = 0.1196581196581197
$ws.Range("M8").Value = 0.101
$ws.Range("N8").Value = 0.0002845872076641308
$ws.Range("O8").Value = 0.003435374149659865
$ws.Range("P8").Value = 0.05
$ws.Range("Q8").Value = 0.0001408847562693717
$ws.Range("R8").Value = 0.001700680272108844
$ws.Range("S8").Value = 0.051
$ws.Range("T8").Value = 0.504950495049505
$ws.Range("U8").Value = 43.5
$ws.Range("V8").Value = 0.1225697379543533
$ws.Range("W8").Value = 0.06480052898391007
$ws.Range("X8").Value = 0.0394118216358498
$ws.Range("Y8").Value = 0.02538870734806028
$ws.Range("Z8").Value = 0.1484771573604061

# Clear removed cells
$ws.Range("F2").ClearContents()
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()
$ws.Range("AQ7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("AN8").ClearContents()
$ws.Range("AO8").ClearContents()
$ws.Range("AP8").ClearContents()
$ws.Range("F8").ClearContents()
